# Protocol change (payment request signature removed)
#
# 1) Bump the datetimeFigureOut placeholder text (master + all 11 layouts)
#    from 2019-09-09 to 2019-09-14.
# 2) Slide 1: the "Signed Payment Request" badge loses its italic/accent6
#    "Signed " run, keeps plain "Payment Request" (now split "Payment "/
#    "Request"), and the shape is repositioned/resized (now centered on the
#    remaining text).
# 3) Slide 1: bump the footer version stamp from V0.97 to V0.98 and its
#    date stamp from 2019-09-09 to 2019-09-14 (splitting the trailing
#    ", 2019-09-09" run into ", " + "2019-09-14").

function EmuToPt($emu) {
    # +0.5 EMU half-step nudge so the runtime's pt->EMU round-trip lands
    # back on the exact integer EMU value instead of flooring one short.
    return ($emu / 12700) + (0.5 / 12700)
}

function Set-DatePlaceholders($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Name -like "Date Placeholder*") {
            $sh.TextFrame.TextRange.Text = "2019-09-14"
        }
    }
}

$p = $ppt.ActivePresentation

# --- 1) Slide master + every layout's date placeholder -------------------
$master = $p.SlideMaster
Set-DatePlaceholders $master.Shapes
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Set-DatePlaceholders $layout.Shapes
}

# --- Slide 1 ---------------------------------------------------------------
$s = $p.Slides.Item(1)

# --- 2) "Signed Payment Request" badge -> "Payment Request" --------------
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.Name -eq "TextBox 191") {
        $sh.TextFrame.TextRange.Text = "Payment Request"
        $tr = $sh.TextFrame.TextRange
        # Re-split into "Payment " / "Request" runs (matches the edited file)
        $firstPart = $tr.Characters(1, 8)
        $firstPart.Text = "Payment "

        $sh.Left = EmuToPt(3062262)
        $sh.Top = EmuToPt(2564904)
        $sh.Width = EmuToPt(1599219)
        $sh.Height = EmuToPt(252948)
    }
}

# --- 3) Footer version/date stamp ------------------------------------------
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.Name -eq "TextBox 279") {
        $tr = $sh.TextFrame.TextRange

        $verRun = $tr.Characters(1, 7)
        $verRun.Text = "V0.98, "

        $commaRun = $tr.Characters(18, 2)
        $commaRun.Text = ", "

        $dateRun = $tr.Characters(20, 10)
        $dateRun.Text = "2019-09-14"
    }
}
